# Update "EDCR Results" sheet with new spike-method result columns (B:H)
# for rows 2-100, per commit "feat: added old spike method results in
# commodity_test_20".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-52 (epsilon 0.001 .. 0.051): pre/recall/F1 updated, NSC/PSC/NRC/PRC
# settle to the "no spike" classification counts.
$ws.Range("B2:B52").Value = 0.1589041095890411
$ws.Range("C2:C52").Value = 1
$ws.Range("D2:D52").Value = 0.2742316784869976
$ws.Range("E2:E52").Value = 0
$ws.Range("F2:F52").Value = 336
$ws.Range("G2:G52").Value = 0
$ws.Range("H2:H52").Value = 5

# Rows 53-100 (epsilon 0.052 .. 0.099): same pre/recall/F1, but NSC/PSC/PRC
# reflect the higher epsilon threshold counts.
$ws.Range("B53:B100").Value = 0.1589041095890411
$ws.Range("C53:C100").Value = 1
$ws.Range("D53:D100").Value = 0.2742316784869976
$ws.Range("E53:E100").Value = 29
$ws.Range("F53:F100").Value = 365
$ws.Range("G53:G100").Value = 5
$ws.Range("H53:H100").Value = 5
